# Update Chocobo_Profits market-data sheets (currentAveragePrice / Leve price / profit columns)
# with refreshed values, per scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1277256.6
$ws.Range("I11").Value = 1277256.6
$ws.Range("K11").Value = 1277256.6
$ws.Range("M11").Value = -1277116.6

$ws.Range("H17").Value = 1451.4108
$ws.Range("J17").Value = 1496.0233
$ws.Range("L17").Value = 4488.0699
$ws.Range("N17").Value = -4824.0699

$ws.Range("H38").Value = 3812.3845
$ws.Range("I38").Value = 123
$ws.Range("J38").Value = 8116.6665
$ws.Range("K38").Value = 369
$ws.Range("L38").Value = 24349.9995
$ws.Range("M38").Value = 3
$ws.Range("N38").Value = -25093.9995

$ws.Range("H39").Value = 518.4
$ws.Range("I39").Value = 184.5
$ws.Range("J39").Value = 900
$ws.Range("K39").Value = 553.5
$ws.Range("L39").Value = 2700
$ws.Range("M39").Value = -257.5
$ws.Range("N39").Value = -3292

$ws.Range("H113").Value = 13857.857
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 13857.857
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 13857.857
$ws.Range("M113").ClearContents() | Out-Null
$ws.Range("N113").Value = -20365.857

$ws.Range("H117").Value = 38896.668
$ws.Range("J117").Value = 38896.668
$ws.Range("L117").Value = 38896.668
$ws.Range("N117").Value = -48074.668

$ws.Range("H132").Value = 83908.08
$ws.Range("I132").Value = 88641.22
$ws.Range("K132").Value = 265923.66
$ws.Range("M132").Value = -263393.66

$ws.Range("H137").Value = 2775.7646
$ws.Range("I137").Value = 2011.5416
$ws.Range("J137").Value = 4609.9
$ws.Range("K137").Value = 6034.6248
$ws.Range("L137").Value = 13829.7
$ws.Range("M137").Value = -3484.6248
$ws.Range("N137").Value = -18929.7


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7993.349
$ws.Range("I32").Value = 5268.2285
$ws.Range("J32").Value = 11399.75
$ws.Range("K32").Value = 5268.2285
$ws.Range("L32").Value = 11399.75
$ws.Range("M32").Value = -4981.2285
$ws.Range("N32").Value = -11973.75

$ws.Range("H74").Value = 7858.875
$ws.Range("I74").Value = 9908.200000000001
$ws.Range("J74").Value = 4443.3335
$ws.Range("K74").Value = 9908.200000000001
$ws.Range("L74").Value = 4443.3335
$ws.Range("M74").Value = -9034.200000000001
$ws.Range("N74").Value = -6191.3335

$ws.Range("H77").Value = 7858.875
$ws.Range("I77").Value = 9908.200000000001
$ws.Range("J77").Value = 4443.3335
$ws.Range("K77").Value = 49541
$ws.Range("L77").Value = 22216.6675
$ws.Range("M77").Value = -45173
$ws.Range("N77").Value = -30952.6675


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 35719616
$ws.Range("I31").Value = 2200
$ws.Range("J31").Value = 41672516
$ws.Range("K31").Value = 2200
$ws.Range("L31").Value = 41672516
$ws.Range("M31").Value = -1905
$ws.Range("N31").Value = -41673106

$ws.Range("H34").Value = 35719616
$ws.Range("I34").Value = 2200
$ws.Range("J34").Value = 41672516
$ws.Range("K34").Value = 2200
$ws.Range("L34").Value = 41672516
$ws.Range("M34").Value = -1998
$ws.Range("N34").Value = -41672920

$ws.Range("H58").Value = 1808.3507
$ws.Range("I58").Value = 1640.541
$ws.Range("J58").Value = 2448.125
$ws.Range("K58").Value = 1640.541
$ws.Range("L58").Value = 2448.125
$ws.Range("M58").Value = -1437.541
$ws.Range("N58").Value = -2854.125

$ws.Range("H82").Value = 39700
$ws.Range("J82").Value = 39700
$ws.Range("L82").Value = 39700
$ws.Range("N82").Value = -40422

$ws.Range("H85").Value = 39700
$ws.Range("J85").Value = 39700
$ws.Range("L85").Value = 39700
$ws.Range("N85").Value = -42196

$ws.Range("H122").Value = 1945.6364
$ws.Range("I122").Value = 1199.3334
$ws.Range("J122").Value = 3544.8572
$ws.Range("K122").Value = 3598.0002
$ws.Range("L122").Value = 10634.5716
$ws.Range("M122").Value = -1148.0002
$ws.Range("N122").Value = -15534.5716

$ws.Range("H135").Value = 38945
$ws.Range("J135").Value = 38945
$ws.Range("L135").Value = 38945
$ws.Range("N135").Value = -49085

$ws.Range("H136").Value = 1808.3507
$ws.Range("I136").Value = 1640.541
$ws.Range("J136").Value = 2448.125
$ws.Range("K136").Value = 4921.623
$ws.Range("L136").Value = 7344.375
$ws.Range("M136").Value = -2371.623
$ws.Range("N136").Value = -12444.375


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1671
$ws.Range("I5").Value = 266
$ws.Range("J5").Value = 4225.5454
$ws.Range("K5").Value = 798
$ws.Range("L5").Value = 12676.6362
$ws.Range("M5").Value = -686
$ws.Range("N5").Value = -12900.6362

$ws.Range("H124").Value = 7666.6665
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 7666.6665
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 22999.9995
$ws.Range("M124").ClearContents() | Out-Null
$ws.Range("N124").Value = -32819.99950000001

$ws.Range("H131").Value = 973.6774
$ws.Range("I131").Value = 2812.5
$ws.Range("J131").Value = 846.86206
$ws.Range("K131").Value = 8437.5
$ws.Range("L131").Value = 2540.58618
$ws.Range("M131").Value = -3397.5
$ws.Range("N131").Value = -12620.58618

$ws.Range("H135").Value = 1671
$ws.Range("I135").Value = 266
$ws.Range("J135").Value = 4225.5454
$ws.Range("K135").Value = 2394
$ws.Range("L135").Value = 38029.9086
$ws.Range("M135").Value = 141
$ws.Range("N135").Value = -43099.9086


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 44674.363
$ws.Range("J133").Value = 44674.363
$ws.Range("L133").Value = 44674.363
$ws.Range("N133").Value = -54794.363


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1685.375
$ws.Range("I46").Value = 1407.238
$ws.Range("K46").Value = 1407.238
$ws.Range("M46").Value = -1219.238

$ws.Range("H81").Value = 61075.6
$ws.Range("J81").Value = 61075.6
$ws.Range("L81").Value = 61075.6
$ws.Range("N81").Value = -63071.6

$ws.Range("H84").Value = 61075.6
$ws.Range("J84").Value = 61075.6
$ws.Range("L84").Value = 183226.8
$ws.Range("N84").Value = -193210.8

$ws.Range("H132").Value = 4756.306
$ws.Range("I132").Value = 1906.25
$ws.Range("J132").Value = 6138.1514
$ws.Range("K132").Value = 5718.75
$ws.Range("L132").Value = 18414.4542
$ws.Range("M132").Value = -3188.75
$ws.Range("N132").Value = -23474.4542


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 14851
$ws.Range("J56").Value = 22504.8
$ws.Range("L56").Value = 22504.8
$ws.Range("N56").Value = -23932.8

$ws.Range("H130").Value = 38165.6
$ws.Range("J130").Value = 38165.6
$ws.Range("L130").Value = 38165.6
$ws.Range("N130").Value = -48205.6

$ws.Range("H133").Value = 60500
$ws.Range("J133").Value = 60500
$ws.Range("L133").Value = 60500
$ws.Range("N133").Value = -70620

